$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet title (workbook.xml)
$ws.Name = "Through 2021-11-24"

# Update label for November row (A12) to reflect new cutoff date
$ws.Range("A12").Value = "November (through 11-24)"

# Update November row (row 12) values
$ws.Range("B12").Value = 24
$ws.Range("C12").Value = 60
$ws.Range("D12").Value = 91
$ws.Range("E12").Value = 48
$ws.Range("F12").Value = 43
$ws.Range("G12").Value = 174
$ws.Range("H12").Value = 165

# Update Total row (row 13) values
$ws.Range("B13").Value = 282
$ws.Range("C13").Value = 546
$ws.Range("D13").Value = 801
$ws.Range("E13").Value = 663
$ws.Range("F13").Value = 525
$ws.Range("G13").Value = 1231
$ws.Range("H13").Value = 1608
